$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.872.62"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").Value = "'211.06"
$ws.Range("E5").Value = "  -1.17%  "
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'23.43"
$ws.Range("E8").Value = "  -0.82%  "
$ws.Range("E9").Value = "  -2.21%  "
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").Value = "1.856.07"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("D13").Value = "1.621.95"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("E14").Value = "  -1.83%  "
$ws.Range("E15").Value = "  -2.30%  "
$ws.Range("D16").Value = "'65.43"
$ws.Range("E16").Value = "  -0.76%  "
$ws.Range("D17").Value = "27.859.59"
$ws.Range("E17").Value = "  -0.44%  "
$ws.Range("D18").Value = "'230.07"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("D20").Value = "0.0₃0722"
$ws.Range("E20").Value = "  -0.26%  "
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").Value = "'10.16"
$ws.Range("E23").Value = "  -5.66%  "
$ws.Range("E24").Value = "  -2.65%  "
$ws.Range("D25").Value = "'154.61"
$ws.Range("E25").Value = "  +1.98%  "
$ws.Range("E26").Value = "  -1.10%  "
$ws.Range("E27").Value = "  -0.20%  "
$ws.Range("D28").Value = "'15.54"
$ws.Range("E28").Value = "  -1.22%  "
$ws.Range("E29").Value = "  -0.25%  "
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("D34").Value = "1.398.59"
$ws.Range("E34").Value = "  -1.22%  "
$ws.Range("D35").Value = "'1.57"
$ws.Range("E35").Value = "  -0.41%  "
$ws.Range("D36").Value = "'0.997"
$ws.Range("E36").Value = "  +8.95%  "
$ws.Range("E37").Value = "  -1.23%  "
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").Value = "'0.555"
$ws.Range("E39").Value = "  -0.49%  "
$ws.Range("D40").Value = "'0.861"
$ws.Range("E40").Value = "  -2.94%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("E43").Value = "  -2.04%  "
$ws.Range("E44").Value = "  +0.13%  "
$ws.Range("D45").Value = "'65.86"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").Value = "1.766.66"
$ws.Range("E47").Value = "  -0.90%  "
$ws.Range("D48").Value = "'88.03"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("E49").Value = "  +1.52%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "'0.0503"
$ws.Range("E50").Value = "  -0.52%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.53"
$ws.Range("E51").Value = "  -1.13%  "
